$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the existing sheet from "Sheet1" to "Data" ---
$ws.Name = "Data"

# --- Update the "Data" sheet contents (row order + new values) ---
# Set columns in the order: A (paths), B (code blocks), D (bare "@Test"),
# then C (plain numbers) so the shared-string table is built up in the
# same sequence as the target workbook.

$ws.Range("A2").Value = "C:\Users\320074769\Downloads\My_Repos\Python_Repos\functiondefextractor\test_resource\test_repo\src\CerberusTest.java_testCerebruswithOutArguments"
$ws.Range("A3").Value = "C:\Users\320074769\Downloads\My_Repos\Python_Repos\functiondefextractor\test_resource\test_repo\src\CerberusTest.java_testCerebrusWithArguments"
$ws.Range("A4").Value = "C:\Users\320074769\Downloads\My_Repos\Python_Repos\functiondefextractor\test_resource\test_repo\src\CerberusTest.java_testCerebruswithWrongArguments"
$ws.Range("A5").Value = "C:\Users\320074769\Downloads\My_Repos\Python_Repos\functiondefextractor\test_resource\test_repo\src\CerberusTest.java_testCallMethod"

$ws.Range("B2").Value = "@Test`r`npublic void testCerebruswithOutArguments() {`r`nCerberus.main(new String[] {})`r`nString expectedOutputString = getCerberusCommandLineUsageString()`r`nassertEquals(expectedOutputString, getModifiedOutputStream().toString())`r`n}`r`n"
$ws.Range("B3").Value = "@Test`r`npublic void testCerebrusWithArguments() {`r`ngetOriginalOutputStream().flush()`r`nCerberus.main(new String[] { `"CPD`" })`r`n}`r`n"
$ws.Range("B4").Value = "@Test`r`npublic void testCerebruswithWrongArguments() {`r`nString dummyArgument = `"dummy argument`"`r`nCerberus.main(new String[] { dummyArgument })`r`nString expectedOutputString = new StringBuilder().append(`"Unmatched argument at index 0: 'dummy argument'`").append(NEW_LINE).append(getCerberusCommandLineUsageString()).toString()`r`n}`r`n"
$ws.Range("B5").Value = "@Test`r`npublic void testCallMethod() throws Exception {`r`nassertEquals(Integer.valueOf(0), new Cerberus().call())`r`n}`r`n"

$ws.Range("D2").Value = "@Test`r`n"
$ws.Range("D3").Value = "@Test`r`n"
$ws.Range("D4").Value = "@Test`r`n"
$ws.Range("D5").Value = "@Test`r`n"

$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 1

# --- Add a new worksheet "Pivot Table" right after "Data" ---
$wsPivot = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$wsPivot.Name = "Pivot Table"

$wsPivot.Range("B1").Value = "@Test Statements"
$wsPivot.Range("C1").Value = "Different @Test pattern counts"
$wsPivot.Range("B2").Value = "@Test"
$wsPivot.Range("A2").Value = 0
$wsPivot.Range("C2").Value = 4

# Apply the same bold/centered/bordered header style used on the "Data"
# sheet's header row (A1:D1) to the new header-like cells B1:C1 and A2.
$ws.Range("A1").Copy()
$wsPivot.Range("B1:C1").PasteSpecial(-4122)  # xlPasteFormats
$wsPivot.Range("A2").PasteSpecial(-4122)     # xlPasteFormats

# Make "Data" the active sheet again (matches original workbook view)
$ws.Activate()
